# Insert a new row for "Ghost" above the current row 21 (Hydralisk),
# shifting all subsequent unit rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the "Ghost" unit data.
$ws.Range("A21").Value = "Ghost"

# Formula columns mirror the pattern used by every other data row.
$ws.Range("B21").Formula = "=D21*C21"
$ws.Range("C21").Formula = "=K21/AVERAGE(K2:K65535)"
$ws.Range("D21").Formula = "=L21/AVERAGE(L2:L65535)"
$ws.Range("E21").Formula = "=N21*(T21+U21+M21)"
$ws.Range("F21").Formula = "=P21*(O21+Z21*R21*((AA21+AB21)/2+Q21)/(AC21/S21))"

$ws.Range("G21").Value = 350
$ws.Range("H21").Value = 245
$ws.Range("I21").Value = 7
$ws.Range("J21").Value = 12.25

$ws.Range("K21").Formula = "=E21/G21"
$ws.Range("L21").Formula = "=F21/G21"

$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 1
$ws.Range("O21").Value = 25
$ws.Range("P21").Value = 1
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 1
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 575
$ws.Range("U21").Value = 0
$ws.Range("V21").Value = "Light"
$ws.Range("W21").Value = 5
$ws.Range("X21").Value = "Magic"
$ws.Range("Y21").Value = 4
$ws.Range("Z21").Value = 1
$ws.Range("AA21").Value = 30
$ws.Range("AB21").Value = 33
$ws.Range("AC21").Value = 0.9
$ws.Range("AD21").Value = 2.5
$ws.Range("AE21").Value = "Detonates an EMP at target point that strips 100.00 + 5.00% shields from enemies."

# Ghost has no abil2/abil3 text, matching every other unit's blank ability
# columns.
$ws.Range("AF21").Value = ""
$ws.Range("AG21").Value = ""

Write-Host "Ghost row inserted"
